# ---------------------------------------------------------------------------
# Applies the resume edits described by the commit "tweeked resume and intro":
#   1. "Sr. Software Engineer"  -> "Software Engineer"
#   2. "Create" (first run of the last Avid Technology bullet) -> "Design and develop"
#      (all the other, already-separate runs making up that bullet must stay
#      untouched/unmerged)
#   3. Two new bullet points added after "...without assistance from software
#      developers" and before the following "Environment: ..." bullet:
#        - "Document detailed design specifications and users guide"
#        - "Installation" + ", training and customer support"  (two runs)
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# --- helper: split the run that currently spans [start,end) into two runs by
# --- toggling (and immediately reverting) a character formatting attribute
# --- across the *second* piece. The engine coalesces adjacent same-formatted
# --- runs on every write, so nudging a property on/off again is the only way
# --- to force a run boundary to "stick" once both pieces end up with the same
# --- resolved formatting again.
function Seal-RunBoundary($rangeToSeal) {
    $rangeToSeal.Font.Bold = 1
    $rangeToSeal.Font.Bold = 0
}

# ===========================================================================
# 1) "Sr. Software Engineer" -> "Software Engineer"
#    (unique in the document, single run in its own paragraph -> plain
#    Find/Replace is safe here, nothing else in that paragraph to disturb.)
# ===========================================================================
$d.Content.Find.Execute("Sr. Software Engineer", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Software Engineer", 2) | Out-Null

# ===========================================================================
# 2) "Create" -> "Design and develop"
#    This run is immediately followed by a dozen more runs that make up the
#    rest of that bullet's sentence. Editing the "Create" run's text directly
#    would coalesce everything after it (through the end of the paragraph)
#    into a single run, which the diff does NOT do - every other run in that
#    bullet is untouched. So: replace the text, then re-split the remainder
#    back along its original run boundaries (their text/formatting is
#    unchanged, only their grouping into <w:r> elements needs restoring).
# ===========================================================================
$rng = $d.Content
$found = $rng.Find.Execute("Create")
if ($found) {
    $rng.Text = "Design and develop"
    $afterStart = $rng.End

    # Lengths (in characters) of the original runs that followed "Create" in
    # the source document - their text/formatting is not changing, only
    # their run segmentation needs to be preserved.
    $segLens = @(2, 4, 7, 12, 4, 7, 1, 16, 7, 26, 10, 43)

    $pos = $afterStart
    foreach ($len in $segLens) {
        $seg = $d.Range($pos, $pos + $len)
        Seal-RunBoundary $seg
        $pos = $pos + $len
    }
}

# ===========================================================================
# 3) Insert the two new bullet points after "...without assistance from
#    software developers" (end of the bullet edited above) and before the
#    following "Environment: ..." bullet.
# ===========================================================================
$rng = $d.Content
$found = $rng.Find.Execute("without assistance from software developers")
if ($found) {
    $rng.Collapse(0)

    # --- New bullet: "Document detailed design specifications and users guide"
    $rng.InsertParagraphAfter()
    $rng.Collapse(0)
    $rng.Move(4, 1) | Out-Null
    $rng.InsertAfter("Document detailed design specifications and users guide")

    # --- New bullet: "Installation" + ", training and customer support"
    # Insert the full sentence in one shot (so there is nothing left to merge
    # it with afterwards), then re-split it into the two runs the diff wants
    # by toggling formatting across just the "Installation" prefix.
    $rng.Collapse(0)
    $rng.InsertParagraphAfter()
    $rng.Collapse(0)
    $rng.Move(4, 1) | Out-Null
    $bulletText = "Installation, training and customer support"
    $rng.InsertAfter($bulletText)

    $bulletStart = $rng.End - $bulletText.Length
    $installSeal = $d.Range($bulletStart, $bulletStart + "Installation".Length)
    Seal-RunBoundary $installSeal
}
